# Update PR log from #25 - append new row 5 to the PR log sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = 25
$ws.Range("B5").Value = "Update index.py"
$ws.Range("C5").Value = "riya-morankar"
$ws.Range("D5").Value = "N/A"
$ws.Range("E5").Value = "edit1 to main"

# The Date column holds a plain "YYYY-MM-DD" text label (matches rows 2-4),
# not an actual date value, so force text formatting before assigning it --
# otherwise Excel auto-parses the literal as a date serial number. Reset the
# cell style back to Normal afterwards so no stray number format lingers.
$ws.Range("F5").NumberFormat = "@"
$ws.Range("F5").Value = "2025-06-17"
$ws.Range("F5").Style = "Normal"
